$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ICABY")

$ws.Range("D8").Value = 160300
$ws.Range("E8").Value = 179200
$ws.Range("F8").Value = 192300
$ws.Range("G8").Value = 212200
$ws.Range("H8").Value = 246100
$ws.Range("I8").Value = 271000
$ws.Range("J8").Value = 268700

$ws.Range("D9").Value = 134900
$ws.Range("E9").Value = 134300
$ws.Range("F9").Value = 151100
$ws.Range("G9").Value = 158500
$ws.Range("H9").Value = 179700
$ws.Range("I9").Value = 224200
$ws.Range("J9").Value = 206700

$ws.Range("D10").Value = 25400
$ws.Range("E10").Value = 44800
$ws.Range("F10").Value = 41200
$ws.Range("G10").Value = 53700
$ws.Range("H10").Value = 66400
$ws.Range("I10").Value = 46800
$ws.Range("J10").Value = 62100

$ws.Range("D14").Value = -9100

$ws.Range("D15").Value = 42500
$ws.Range("E15").Value = 41300
$ws.Range("F15").Value = 28400
$ws.Range("G15").Value = 28900
$ws.Range("H15").Value = 30200
$ws.Range("I15").Value = 32200
$ws.Range("J15").Value = 32900

$ws.Range("D17").Value = 207800
$ws.Range("E17").Value = 218900
$ws.Range("F17").Value = 223200
$ws.Range("G17").Value = 229900
$ws.Range("H17").Value = 257300
$ws.Range("I17").Value = 306000
$ws.Range("J17").Value = 291600

$ws.Range("D18").Value = -47400
$ws.Range("E18").Value = -39700
$ws.Range("F18").Value = -30900
$ws.Range("G18").Value = -17700
$ws.Range("H18").Value = -11200
$ws.Range("I18").Value = -35000
$ws.Range("J18").Value = -22900

$ws.Range("F21").Value = 12700
$ws.Range("G21").Value = 26300
$ws.Range("H21").Value = 30400
$ws.Range("J21").Value = 20100

$ws.Range("D23").Value = -46700
$ws.Range("E23").Value = -40500
$ws.Range("F23").Value = -31400
$ws.Range("G23").Value = -17700
$ws.Range("H23").Value = -11500
$ws.Range("I23").Value = -34800
$ws.Range("J23").Value = -23300

$ws.Range("D26").Value = -46200
$ws.Range("E26").Value = -39800
$ws.Range("F26").Value = -29700
$ws.Range("G26").Value = -17800
$ws.Range("H26").Value = -11800
$ws.Range("I26").Value = -35500
$ws.Range("J26").Value = -22900

$ws.Range("D27").Value = -46200
$ws.Range("E27").Value = -39800
$ws.Range("F27").Value = -29700
$ws.Range("G27").Value = -17800
$ws.Range("H27").Value = -11800
$ws.Range("I27").Value = -35500
$ws.Range("J27").Value = -22900

$ws.Range("D33").Value = -46200
$ws.Range("E33").Value = -39800
$ws.Range("F33").Value = -29700
$ws.Range("G33").Value = -17800
$ws.Range("H33").Value = -11800
$ws.Range("I33").Value = -35500
$ws.Range("J33").Value = -22900

$ws.Range("D35").Value = -46200
$ws.Range("E35").Value = -39800
$ws.Range("F35").Value = -29700
$ws.Range("G35").Value = -17800
$ws.Range("H35").Value = -11800
$ws.Range("I35").Value = -35500
$ws.Range("J35").Value = -22900

$ws.Range("D41").Value = 74500
$ws.Range("F41").Value = 10500
$ws.Range("G41").Value = 7900
$ws.Range("H41").Value = 23200
$ws.Range("I41").Value = 36700
$ws.Range("J41").Value = 43100

$ws.Range("G43").Value = 10000

$ws.Range("J44").Value = 700

$ws.Range("E45").Value = 4400
$ws.Range("G45").Value = 10600
$ws.Range("I45").Value = 7900
$ws.Range("J45").Value = 20600

$ws.Range("D46").Value = 86700
$ws.Range("E46").Value = 26800
$ws.Range("F46").Value = 24700
$ws.Range("G46").Value = 30900
$ws.Range("H46").Value = 42200
$ws.Range("I46").Value = 57000
$ws.Range("J46").Value = 75700

$ws.Range("D48").Value = 113200
$ws.Range("E48").Value = 113900
$ws.Range("F48").Value = 111300
$ws.Range("G48").Value = 113800
$ws.Range("H48").Value = 119000
$ws.Range("I48").Value = 130500
$ws.Range("J48").Value = 136900

$ws.Range("D49").Value = 19500
$ws.Range("E49").Value = 21800
$ws.Range("F49").Value = 20400
$ws.Range("G49").Value = 21900
$ws.Range("H49").Value = 18000
$ws.Range("J49").Value = 14100

$ws.Range("D52").Value = 43600
$ws.Range("E52").Value = 44600
$ws.Range("F52").Value = 47300
$ws.Range("G52").Value = 46000
$ws.Range("H52").Value = 46800
$ws.Range("I52").Value = 48600
$ws.Range("J52").Value = 50300

$ws.Range("D54").Value = 263000
$ws.Range("E54").Value = 207100
$ws.Range("F54").Value = 203700
$ws.Range("G54").Value = 212600
$ws.Range("H54").Value = 226000
$ws.Range("I54").Value = 250400
$ws.Range("J54").Value = 277000

$ws.Range("E57").Value = 7500
$ws.Range("H57").Value = 11400
$ws.Range("J57").Value = 11900

$ws.Range("D58").Value = 50300
$ws.Range("E58").Value = 75200
$ws.Range("F58").Value = 38200
$ws.Range("G58").Value = 12700

$ws.Range("D59").Value = 50700
$ws.Range("E59").Value = 57700
$ws.Range("F59").Value = 53600
$ws.Range("G59").Value = 52500
$ws.Range("H59").Value = 57300
$ws.Range("I59").Value = 53700
$ws.Range("J59").Value = 58000

$ws.Range("D60").Value = 103200
$ws.Range("E60").Value = 140400
$ws.Range("F60").Value = 96100
$ws.Range("G60").Value = 74100
$ws.Range("H60").Value = 68700
$ws.Range("I60").Value = 67300
$ws.Range("J60").Value = 69900

$ws.Range("D61").Value = 12700
$ws.Range("I61").Value = 12700

$ws.Range("H62").Value = 5400
$ws.Range("J62").Value = 8000

$ws.Range("D66").Value = 119100
$ws.Range("E66").Value = 142700
$ws.Range("F66").Value = 99400
$ws.Range("G66").Value = 78500
$ws.Range("H66").Value = 74100
$ws.Range("I66").Value = 86700
$ws.Range("J66").Value = 77900

$ws.Range("D72").Value = -855800
$ws.Range("E72").Value = -809600
$ws.Range("F72").Value = -769700
$ws.Range("G72").Value = -740000
$ws.Range("H72").Value = -721300
$ws.Range("I72").Value = -709400
$ws.Range("J72").Value = -674000

$ws.Range("D76").Value = 143900
$ws.Range("E76").Value = 64400
$ws.Range("F76").Value = 104300
$ws.Range("G76").Value = 134100
$ws.Range("H76").Value = 151900
$ws.Range("I76").Value = 163600
$ws.Range("J76").Value = 199100

$ws.Range("D81").Value = -46200
$ws.Range("E81").Value = -39800
$ws.Range("F81").Value = -29700
$ws.Range("G81").Value = -17800
$ws.Range("H81").Value = -11800
$ws.Range("I81").Value = -35500
$ws.Range("J81").Value = -22900

$ws.Range("D83").Value = 41300
$ws.Range("E83").Value = 44100
$ws.Range("F83").Value = 44000
$ws.Range("G83").Value = 41800
$ws.Range("H83").Value = 44100
$ws.Range("I83").Value = 43400
$ws.Range("J83").Value = "NA"

$ws.Range("D89").Value = 11400
$ws.Range("E89").Value = 17400
$ws.Range("F89").Value = 14500
$ws.Range("G89").Value = 32600
$ws.Range("H89").Value = 20600
$ws.Range("I89").Value = 24600
$ws.Range("J89").Value = 30900

$ws.Range("D91").Value = -29800
$ws.Range("E91").Value = -25100
$ws.Range("F91").Value = -23900
$ws.Range("G91").Value = -18800
$ws.Range("H91").Value = -27300
$ws.Range("I91").Value = -28400
$ws.Range("J91").Value = -30600

$ws.Range("D94").Value = -47200
$ws.Range("E94").Value = -40300
$ws.Range("F94").Value = -38700
$ws.Range("G94").Value = -37200
$ws.Range("H94").Value = -39800
$ws.Range("I94").Value = -38000
$ws.Range("J94").Value = "NA"

$ws.Range("D100").Value = 36200
$ws.Range("E100").Value = 25500
$ws.Range("F100").Value = 12700
$ws.Range("G100").Value = -12700
$ws.Range("H100").Value = 12700
$ws.Range("I100").Value = -500
$ws.Range("J100").Value = "NA"

$ws.Range("I101").Value = 100
$ws.Range("J101").Value = "NA"

$ws.Range("D102").Value = 400
$ws.Range("E102").Value = 2600
$ws.Range("F102").Value = -11400
$ws.Range("G102").Value = -17400
$ws.Range("H102").Value = -6400
$ws.Range("I102").Value = -13800
$ws.Range("J102").Value = -10700
